$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" column (K) was added to the table, mirroring the formatting
# already used by the 2019 column (J) for each populated row.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial(-4122)

$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Fill in the new values for 2020.
$ws.Range("K4").Value = 2020
$ws.Range("K6").Value = 5.9
$ws.Range("K7").Value = 1.5
$ws.Range("K8").Value = "-"

# Restore the selection that was active when the workbook was last saved.
[void]$ws.Range("L16").Select()
